# Edit: rename headers on existing sheets and add new "PO Forecast" sheet
# with weekly PO forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

# 1) Rename header labels on the existing sheets.
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add a new worksheet "PO Forecast" after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse existing cell formatting (bold/border/centered header style, and the
# date number-format used for column A) by copying formats from the
# "Weekly Quantity" sheet instead of constructing brand-new styles.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$forecastRows = @(
    @(44955.99999999999,32,-10.85993080341517,75.55758078583801),
    @(44962.99999999999,33,-12.84977127189371,79.74256002576003),
    @(44969.99999999999,33,-12.81090569447243,73.1419728709415),
    @(44976.99999999999,33,-8.929748763774786,74.91533224822346),
    @(44983.99999999999,34,-9.724710080521538,78.27238762422658),
    @(44990.99999999999,34,-7.620933516224,77.20854660422394),
    @(44997.99999999999,34,-10.3538272202614,78.5806503796362),
    @(45004.99999999999,35,-7.71178063774533,78.27686856442385),
    @(45011.99999999999,35,-7.561799194461335,80.38585345170493),
    @(45032.99999999999,36,-3.567434336462332,79.6699888027627),
    @(45039.99999999999,37,-4.368682203417563,81.5300247288432),
    @(45053.99999999999,38,-4.937601849547096,82.11363607061691),
    @(45060.99999999999,38,0.002199687694895411,85.0573743549179),
    @(45067.99999999999,38,-0.6764043573730033,84.10265580962587),
    @(45088.99999999999,39,-5.985603065446033,82.87305689059856),
    @(45095.99999999999,40,-5.113206063348859,81.07136750699448),
    @(45102.99999999999,40,-2.779884957052984,82.85246780761285),
    @(45109.99999999999,41,-5.88608566600058,80.83577981639884),
    @(45144.99999999999,43,0.8528602271784486,85.18614683784664),
    @(45158.99999999999,43,1.28555925744528,89.76776763807382),
    @(45165.99999999999,44,-3.058074901873497,90.65092023628563),
    @(45179.99999999999,44,2.401392311827806,86.55948294914863),
    @(45186.99999999999,45,-1.295710541775224,87.43269756392857),
    @(45207.99999999999,46,2.758464559031924,91.40358125763645),
    @(45214.99999999999,46,2.7595893329133,90.50245066927579),
    @(45221.99999999999,47,5.814424121967414,92.3212418503917),
    @(45228.99999999999,47,2.359246190567024,91.88624403886784),
    @(45242.99999999999,48,4.084827479256883,90.68153658776866),
    @(45249.99999999999,48,5.505544237046595,94.8902112766622),
    @(45270.99999999999,49,8.618216953874038,96.26890932379128),
    @(45277.99999999999,50,4.506424642975058,93.13870362518158),
    @(45298.99999999999,51,11.48608788281785,94.68419205509542),
    @(45305.99999999999,51,9.009739084841261,95.23739824685281),
    @(45319.99999999999,52,11.89053638624672,94.16805686869651),
    @(45326.99999999999,53,10.82492708397616,95.14163945397446),
    @(45354.99999999999,54,12.33198100357547,97.99835699534147),
    @(45361.99999999999,54,10.74566662397709,98.4525757113005),
    @(45368.99999999999,55,11.82567646755736,97.22027177120496),
    @(45382.99999999999,56,9.841660703444381,102.7369275275635),
    @(45389.99999999999,56,10.56111655268359,100.2363150924824),
    @(45403.99999999999,57,12.58227305073024,98.94094092768762),
    @(45501.99999999999,62,15.85368478269511,103.5541155575336),
    @(45515.99999999999,63,17.6936243939304,105.0715399141462),
    @(45571.99999999999,66,23.01256738198197,107.9769106004481),
    @(45578.99999999999,66,17.68730420178336,107.6144365668098),
    @(45585.99999999999,67,22.37326728011341,111.494431584454),
    @(45592.99999999999,67,22.18761511900923,108.1519067135563),
    @(45599.99999999999,68,25.93082255245277,110.3540755878035),
    @(45606.99999999999,68,22.4231143953041,112.5114234835408),
    @(45613.99999999999,68,26.24803020983996,113.1547316185315),
    @(45620.99999999999,69,25.19587666475322,108.99151345695),
    @(45627.99999999999,69,26.1725493200844,116.2594920599615)
)

$r = 2
foreach ($row in $forecastRows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r += 1
}
